$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data (A1:B8) so leftover cells from rows 5-8 are removed
$ws.Range("A1:B8").Clear()

# Insert a new column B for "User" (shifts old Status column from B to C)
$ws.Range("B1").EntireColumn.Insert()

# Header row
$ws.Range("A1").Value = "CommitteeName"
$ws.Range("B1").Value = "User"
$ws.Range("C1").Value = "Status"

# Data rows (shared-string table must grow in the same order as the target file,
# so set B3 before B2)
$ws.Range("A2").Value = "aa"
$ws.Range("A3").Value = "bb"
$ws.Range("B3").Value = "ketan Sali,aaaaaa"
$ws.Range("B2").Value = "Harshita sharma,Akshay Bhagwat"
$ws.Range("C2").Value = "Active"
$ws.Range("C3").Value = "DeActive"

$ws.Range("A4").Value = "cc"
$ws.Range("B4").Value = "Himanshu,juku"
$ws.Range("C4").Value = "Active"

# Autofit column A to match bestFit width
$ws.Range("A1:A4").EntireColumn.AutoFit()

# Update selection to match target state
$ws.Range("C15").Select()

$wb.Save()
